$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()

$ws.Range("C5:D5").ClearFormats()

$ws.Range("A5").Value = 102
$ws.Range("B5").Value = "LC/GFG"
$ws.Range("C5").Value = "Binary Tree Level Order Traversal"
$ws.Range("D5").Value = "Java/Python"

$ws.Range("D5").Select()
